# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# Stamp the M2Doc version that generated/validated this template into
# the document's custom properties, the same way Word records things
# like "Company" or "Manager" under Document.CustomDocumentProperties.

$d = $word.ActiveDocument

$propertyName = "M2DocVersion"
$m2docVersion = "3.5.0"

# msoPropertyTypeString = 4
$msoPropertyTypeString = 4

$customProps = $d.CustomDocumentProperties

# Keep this idempotent: drop any previous value for this property
# before (re)adding it, in case the template was already stamped.
ForEach ($existing in $customProps) {
    if ($existing.Name -eq $propertyName) {
        $existing.Delete()
    }
}

$customProps.Add($propertyName, $false, $msoPropertyTypeString, $m2docVersion)

Write-Output "M2Doc version '$m2docVersion' stored in custom property '$propertyName'."
